# Insert a new data row for the Macroferia Regional de Talca - Mango sheet.
# The new record is placed right after the existing row 29, pushing the
# former rows 30..110 down to 31..111 (plain Excel row insert handles the
# shifting of every column/value/format automatically).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(30).Insert()

$ws.Range("A30").Value = 5
$ws.Range("B30").Value = "Macroferia Regional de Talca"
$ws.Range("C30").Value = "Maule"
$ws.Range("D30").Value = 44608
$ws.Range("E30").Value = 7
$ws.Range("F30").Value = "Fruta"
$ws.Range("G30").Value = 100108
$ws.Range("H30").Value = "Tropicales y subtropicales"
$ws.Range("I30").Value = 100108002
$ws.Range("J30").Value = "Mango"
$ws.Range("K30").Value = "Sin especificar"
$ws.Range("L30").Value = "Primera"
$ws.Range("M30").Value = 300
$ws.Range("N30").Value = 7000
$ws.Range("O30").Value = 7000
$ws.Range("P30").Value = 7000
$ws.Range("Q30").Value = "`$/bandeja 4 kilos"
$ws.Range("R30").Value = "Perú"
$ws.Range("S30").Value = 1750
$ws.Range("T30").Value = 4
